$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Brn_Code column (C) values from text "195" to numeric 19
$ws.Range("C2").Value = 19
$ws.Range("C3").Value = 19
$ws.Range("C4").Value = 19

# Fill in Out_Standing_Amount (I) and Disbursed_Amount (J) for rows 2-4
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 50000

$ws.Range("I3").Value = 200
$ws.Range("J3").Value = 20000

$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 20000

# Update the active selection to J12
$ws.Range("J12").Select()
